# Refresh the cryptocurrency price/volume snapshot on Sheet1 (GitHub Actions bot update).
# Column D = Price, Column E = Volume(1h) change. Cells that hold a plain
# decimal-looking price (e.g. "211.39") are first forced to Text format so
# Excel's auto-type-detection doesn't coerce them into floating point
# numbers (which would both change the cell type and introduce binary
# rounding noise); the cell style is then reset back to "Normal" so no
# stray formatting is left behind. Values that are not number-like
# (multi-dot prices such as "26.647.18", or the "  +0.26%  " style
# percentages) are plain text already and need no special handling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.647.18"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.598.37"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.822.47"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.599.39"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "26.630.29"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.67%  "
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "1.275.40"
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.619"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  +17.99%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "1.735.16"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "
